$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 17:18, copying formatting from the row above (row 16)
$ws.Rows("17:18").Insert()

# Row 17 - CyVerse (Tucson, Arizona, USA)
$ws.Range("A17").Value = "CyVerse"
$ws.Range("B17").Value = "CyVerse"
$ws.Range("C17").Value = "Tucson"
$ws.Range("D17").Value = "Arizona"
$ws.Range("E17").Value = "United States of America"
$ws.Range("F17").Value = "USA"
$ws.Range("G17").Value = 32.253943
$ws.Range("H17").Value = -110.974114
$ws.Range("I17").Value = "Yes"
$ws.Range("G17").ClearFormats()
$ws.Range("H17").ClearFormats()

# Row 18 - ESIIL (Boulder, Colorado, USA)
$ws.Range("B18").Value = "ESIIL"
$ws.Range("C18").Value = "Boulder"
$ws.Range("A18").Value = "Environmental Data Science Innovation & Inclusion Lab"
$ws.Range("D18").Value = "Colorado"
$ws.Range("E18").Value = "United States of America"
$ws.Range("F18").Value = "USA"
$ws.Range("G18").Value = 40.022557999999997
$ws.Range("H18").Value = -105.250169
$ws.Range("I18").Value = "Yes"
$ws.Range("A18").ClearFormats()
$ws.Range("G18").ClearFormats()
$ws.Range("H18").ClearFormats()

# Widen column A to fit the new, longer names (matches a manual column-width resize)
$ws.Range("A1:A18").ColumnWidth = 64.28515625

# Move the active selection like in the authored workbook
[void]$ws.Range("A24").Select()
